$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1096
$ws.Range("I2").Value = 308.25
$ws.Range("J2").Value = 3196.6667
$ws.Range("K2").Value = 308.25
$ws.Range("L2").Value = 3196.6667
$ws.Range("M2").Value = -195.25
$ws.Range("N2").Value = -3422.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3588
$ws.Range("J17").Value = 3671.5789
$ws.Range("L17").Value = 11014.7367
$ws.Range("N17").Value = -11350.7367

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 2837.6
$ws.Range("I42").Value = 995.8
$ws.Range("J42").Value = 6521.2
$ws.Range("K42").Value = 2987.4
$ws.Range("L42").Value = 19563.6
$ws.Range("M42").Value = -2757.4
$ws.Range("N42").Value = -20023.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 873
$ws.Range("I101").Value = 873
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 2619
$ws.Range("L101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2987.5
$ws.Range("I106").Value = 2987.5
$ws.Range("K106").Value = 2987.5
$ws.Range("M106").Value = -2356.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3893.25
$ws.Range("J112").Value = 3991
$ws.Range("L112").Value = 11973
$ws.Range("N112").Value = -14189

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 5120
$ws.Range("I127").Value = 5120
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 15360
$ws.Range("L127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -10400

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1101.4
$ws.Range("I132").Value = 1101.4
$ws.Range("K132").Value = 3304.2
$ws.Range("M132").Value = -774.2000000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5962.4736
$ws.Range("I32").Value = 6015.9443
$ws.Range("K32").Value = 6015.9443
$ws.Range("M32").Value = -5728.9443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2037.4
$ws.Range("I74").Value = 2297
$ws.Range("J74").Value = 999
$ws.Range("K74").Value = 2297
$ws.Range("L74").Value = 999
$ws.Range("M74").Value = -1423
$ws.Range("N74").Value = -2747

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2037.4
$ws.Range("I77").Value = 2297
$ws.Range("J77").Value = 999
$ws.Range("K77").Value = 11485
$ws.Range("L77").Value = 4995
$ws.Range("M77").Value = -7117
$ws.Range("N77").Value = -13731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1639.6
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 64785.2
$ws.Range("J134").Value = 64785.2
$ws.Range("L134").Value = 64785.2
$ws.Range("N134").Value = -74925.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 596.3333
$ws.Range("I80").Value = 916.0909
$ws.Range("J80").Value = 244.6
$ws.Range("K80").Value = 916.0909
$ws.Range("L80").Value = 244.6
$ws.Range("M80").Value = 81.90909999999997
$ws.Range("N80").Value = -2240.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 596.3333
$ws.Range("I83").Value = 916.0909
$ws.Range("J83").Value = 244.6
$ws.Range("K83").Value = 4580.4545
$ws.Range("L83").Value = 1223
$ws.Range("M83").Value = 411.5455000000002
$ws.Range("N83").Value = -11207

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 70777.5
$ws.Range("J122").Value = 70777.5
$ws.Range("L122").Value = 70777.5
$ws.Range("N122").Value = -80577.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 116.26923
$ws.Range("I7").Value = 116.333336
$ws.Range("J7").Value = 116.181816
$ws.Range("K7").Value = 116.333336
$ws.Range("L7").Value = 116.181816
$ws.Range("M7").Value = -3.333336000000003
$ws.Range("N7").Value = -342.181816

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1534.1428
$ws.Range("I31").Value = 1243.75
$ws.Range("K31").Value = 1243.75
$ws.Range("M31").Value = -948.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1534.1428
$ws.Range("I34").Value = 1243.75
$ws.Range("K34").Value = 1243.75
$ws.Range("M34").Value = -1041.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 89834.234
$ws.Range("J94").Value = 4617.4287
$ws.Range("L94").Value = 4617.4287
$ws.Range("N94").Value = -5519.4287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4621.75
$ws.Range("I99").Value = 4499.3335
$ws.Range("K99").Value = 4499.3335
$ws.Range("M99").Value = -3001.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4621.75
$ws.Range("I126").Value = 4499.3335
$ws.Range("K126").Value = 13498.0005
$ws.Range("M126").Value = -11028.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2077.2222
$ws.Range("I132").Value = 2062.5
$ws.Range("K132").Value = 6187.5
$ws.Range("M132").Value = -3657.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2000.375
$ws.Range("I134").Value = 2317.25
$ws.Range("K134").Value = 6951.75
$ws.Range("M134").Value = -4416.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27501468
$ws.Range("I4").Value = 27501468
$ws.Range("K4").Value = 82504404
$ws.Range("M4").Value = -82504292

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 63992
$ws.Range("J37").Value = 63992
$ws.Range("L37").Value = 191976
$ws.Range("N37").Value = -192200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 143.58333
$ws.Range("J40").Value = 396.66666
$ws.Range("L40").Value = 1586.66664
$ws.Range("N40").Value = -1724.66664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 412.5
$ws.Range("I99").Value = 625
$ws.Range("J99").Value = 200
$ws.Range("K99").Value = 1875
$ws.Range("L99").Value = 600
$ws.Range("M99").Value = 371
$ws.Range("N99").Value = -5092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1487.4445
$ws.Range("I113").Value = 1100
$ws.Range("J113").Value = 1797.4
$ws.Range("K113").Value = 3300
$ws.Range("L113").Value = 5392.200000000001
$ws.Range("M113").Value = -1130
$ws.Range("N113").Value = -9732.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 15269.75
$ws.Range("I121").Value = 40486
$ws.Range("J121").Value = 6864.3335
$ws.Range("K121").Value = 121458
$ws.Range("L121").Value = 20593.0005
$ws.Range("M121").Value = -120148
$ws.Range("N121").Value = -23213.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 4493
$ws.Range("I126").Value = 4493
$ws.Range("K126").Value = 13479
$ws.Range("M126").Value = -8539

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 3175.5
$ws.Range("J127").Value = 3175.5
$ws.Range("L127").Value = 9526.5
$ws.Range("N127").Value = -19446.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 16249.5
$ws.Range("I141").Value = 16249.5
$ws.Range("K141").Value = 48748.5
$ws.Range("M141").Value = -43568.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8144.6924
$ws.Range("I70").Value = 7000.125
$ws.Range("J70").Value = 9976
$ws.Range("K70").Value = 7000.125
$ws.Range("L70").Value = 9976
$ws.Range("M70").Value = -6730.125
$ws.Range("N70").Value = -10516

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8144.6924
$ws.Range("I73").Value = 7000.125
$ws.Range("J73").Value = 9976
$ws.Range("K73").Value = 7000.125
$ws.Range("L73").Value = 9976
$ws.Range("M73").Value = -6064.125
$ws.Range("N73").Value = -11848

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1994.8
$ws.Range("I102").Value = 1994.8
$ws.Range("K102").Value = 1994.8
$ws.Range("M102").Value = -372.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3999.8
$ws.Range("I126").Value = 3333
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 9999
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -7529
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1866.3334
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5716.7
$ws.Range("I22").Value = 5513
$ws.Range("J22").Value = 6022.25
$ws.Range("K22").Value = 5513
$ws.Range("L22").Value = 6022.25
$ws.Range("M22").Value = -5218
$ws.Range("N22").Value = -6612.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 5716.7
$ws.Range("I27").Value = 5513
$ws.Range("J27").Value = 6022.25
$ws.Range("K27").Value = 5513
$ws.Range("L27").Value = 6022.25
$ws.Range("M27").Value = -5406
$ws.Range("N27").Value = -6236.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 2999
$ws.Range("I30").Value = 2999
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2999
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -2891

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3759.625
$ws.Range("I40").Value = 2999
$ws.Range("J40").Value = 3868.2856
$ws.Range("K40").Value = 2999
$ws.Range("L40").Value = 3868.2856
$ws.Range("M40").Value = -2863
$ws.Range("N40").Value = -4140.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2897
$ws.Range("I61").Value = 2600
$ws.Range("J61").Value = 3024.2856
$ws.Range("K61").Value = 2600
$ws.Range("L61").Value = 3024.2856
$ws.Range("M61").Value = -2398
$ws.Range("N61").Value = -3428.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2897
$ws.Range("I113").Value = 2600
$ws.Range("J113").Value = 3024.2856
$ws.Range("K113").Value = 2600
$ws.Range("L113").Value = 3024.2856
$ws.Range("M113").Value = -430
$ws.Range("N113").Value = -7364.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 899.125
$ws.Range("I122").Value = 899.6667
$ws.Range("K122").Value = 2699.0001
$ws.Range("M122").Value = -249.0001000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3947.5293
$ws.Range("I126").Value = 2270.9
$ws.Range("K126").Value = 6812.700000000001
$ws.Range("M126").Value = -4342.700000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3597
$ws.Range("I136").Value = 3486.1765
$ws.Range("J136").Value = 3973.8
$ws.Range("K136").Value = 10458.5295
$ws.Range("L136").Value = 11921.4
$ws.Range("M136").Value = -7908.529500000001
$ws.Range("N136").Value = -17021.4
